$d = $word.ActiveDocument

# Fix 1: "den main Branch" -> "dem main Branch" (spelling mistake)
$d.Content.Find.Execute("Der Code auf den main Branch la", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Der Code auf dem main Branch la", 2)

# Fix 2: Replace the incorrect sentence about the Stable-State branch with the corrected one.
$d.Content.Find.Execute(". Nachdem der Code von main am Review auf diesem Branch gesichert wird.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ". Wenn der Code von main beim Review abgenommen wurde, wird dieser in den Stable-State Branch gemerged.", 2)
